$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A42").Value = "Centra"
